$d = $word.ActiveDocument

$old = "1ª persona (               ):"
$new = "1ª persona ( Adrian ):"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
